# Update "想去人数" (want-to-go count) values in column F for the
# "展览" and "全部类型" worksheets, reflecting refreshed counts from the
# generated data source.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 826
$ws1.Range("F6").Value = 140
$ws1.Range("F8").Value = 4969
$ws1.Range("F10").Value = 5232
$ws1.Range("F11").Value = 602
$ws1.Range("F12").Value = 1317

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 826
$ws4.Range("F6").Value = 140
$ws4.Range("F9").Value = 4969
$ws4.Range("F11").Value = 5232
$ws4.Range("F12").Value = 602
$ws4.Range("F13").Value = 1317
